$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C3:C21").NumberFormat = "@"
$ws.Range("C2").Value = "08:00"
$ws.Range("C3").Value = "12:00"
$ws.Range("C4").Value = "15:00"
$ws.Range("C5").Value = "18:00"
